$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7631
$ws1.Range("F13").Value = 3067
$ws1.Range("F15").Value = 82
$ws1.Range("F16").Value = 719
$ws1.Range("F21").Value = 212
$ws1.Range("F22").Value = 214
$ws1.Range("F23").Value = 266
$ws1.Range("F27").Value = 250
$ws1.Range("F32").Value = 28

# Sheet "全部类型" updates (mirrors the same events)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 7631
$ws4.Range("F17").Value = 3067
$ws4.Range("F19").Value = 82
$ws4.Range("F21").Value = 719
$ws4.Range("F27").Value = 212
$ws4.Range("F28").Value = 214
$ws4.Range("F29").Value = 266
$ws4.Range("F33").Value = 250
$ws4.Range("F38").Value = 28
